# Se procesan de nuevo los datos con las nuevas dimensiones curadas
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column F (municipio-nombre) and H (aragon) are now both refArea dimensions,
# just like columns G (provincia-nombre) and J (comarca-nombre) already were.
$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("H2").Value = "sdmx-dimension:refArea"

# F used to be a measure ("medida"); now it's a dimension ("dim") too.
$ws.Range("F3").Value = "dim"

# F (municipio) and H (comunidad/aragon) now reference their own URI columns,
# replacing the previous xsd:int / skos:Concept type markers.
$ws.Range("F4").Value = "URI-Municipio"
$ws.Range("H4").Value = "URI-Comunidad"

# The mapping-aragon.xlsx reference for column H is no longer needed.
$ws.Range("H5").Clear()
